$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4049.8462
$ws.Range("I40").Value = 3449.8572
$ws.Range("K40").Value = 3449.8572
$ws.Range("M40").Value = -3274.8572
$ws.Range("H43").Value = 11730
$ws.Range("I43").Value = 11166.667
$ws.Range("K43").Value = 11166.667
$ws.Range("M43").Value = -11097.667
$ws.Range("H62").Value = 116053.78
$ws.Range("I62").Value = 146784.14
$ws.Range("J62").Value = 8497.5
$ws.Range("K62").Value = 146784.14
$ws.Range("L62").Value = 8497.5
$ws.Range("M62").Value = -146160.14
$ws.Range("N62").Value = -9745.5
$ws.Range("H64").Value = 4023.8823
$ws.Range("I64").Value = 3749.75
$ws.Range("J64").Value = 4681.8
$ws.Range("K64").Value = 3749.75
$ws.Range("L64").Value = 4681.8
$ws.Range("M64").Value = -3501.75
$ws.Range("N64").Value = -5177.8
$ws.Range("H65").Value = 116053.78
$ws.Range("I65").Value = 146784.14
$ws.Range("J65").Value = 8497.5
$ws.Range("K65").Value = 733920.7000000001
$ws.Range("L65").Value = 42487.5
$ws.Range("M65").Value = -730800.7000000001
$ws.Range("N65").Value = -48727.5
$ws.Range("H67").Value = 4023.8823
$ws.Range("I67").Value = 3749.75
$ws.Range("J67").Value = 4681.8
$ws.Range("K67").Value = 3749.75
$ws.Range("L67").Value = 4681.8
$ws.Range("M67").Value = -2891.75
$ws.Range("N67").Value = -6397.8
$ws.Range("H74").Value = 3866.1333
$ws.Range("I74").Value = 2888
$ws.Range("K74").Value = 2888
$ws.Range("M74").Value = -1952
$ws.Range("H77").Value = 3866.1333
$ws.Range("I77").Value = 2888
$ws.Range("K77").Value = 14440
$ws.Range("M77").Value = -9760
$ws.Range("H98").Value = 1063.6154
$ws.Range("I98").Value = 944
$ws.Range("K98").Value = 944
$ws.Range("M98").Value = 554
$ws.Range("H106").Value = 16752.715
$ws.Range("J106").Value = 22398.4
$ws.Range("L106").Value = 22398.4
$ws.Range("N106").Value = -23660.4
$ws.Range("H107").Value = 1485.1428
$ws.Range("I107").Value = 1321
$ws.Range("J107").Value = 1895.5
$ws.Range("K107").Value = 1321
$ws.Range("L107").Value = 1895.5
$ws.Range("M107").Value = 599
$ws.Range("N107").Value = -5735.5
$ws.Range("H112").Value = 3975
$ws.Range("I112").Value = 3512
$ws.Range("J112").Value = 4052.1667
$ws.Range("K112").Value = 10536
$ws.Range("L112").Value = 12156.5001
$ws.Range("M112").Value = -9428
$ws.Range("N112").Value = -14372.5001
$ws.Range("H116").Value = 4747.154
$ws.Range("I116").Value = 4535.6665
$ws.Range("K116").Value = 4535.6665
$ws.Range("M116").Value = -1093.6665
$ws.Range("H122").Value = 1063.6154
$ws.Range("I122").Value = 944
$ws.Range("K122").Value = 2832
$ws.Range("M122").Value = -382
$ws.Range("H132").Value = 2705
$ws.Range("I132").Value = 1469.2858
$ws.Range("K132").Value = 4407.857400000001
$ws.Range("M132").Value = -1877.857400000001
$ws.Range("H137").Value = 4297.25
$ws.Range("I137").Value = 1566.7667
$ws.Range("K137").Value = 4700.300099999999
$ws.Range("M137").Value = -2150.300099999999
$ws.Range("H138").Value = 2705.5293
$ws.Range("J138").Value = 3615.842
$ws.Range("L138").Value = 10847.526
$ws.Range("N138").Value = -21127.526
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 39495
$ws.Range("J54").Value = 39495
$ws.Range("L54").Value = 39495
$ws.Range("N54").Value = -41033
$ws.Range("H61").Value = 985.0606
$ws.Range("I61").Value = 801
$ws.Range("K61").Value = 801
$ws.Range("M61").Value = -589
$ws.Range("H136").Value = 985.0606
$ws.Range("I136").Value = 801
$ws.Range("K136").Value = 2403
$ws.Range("M136").Value = 147
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2726.6667
$ws.Range("I99").Value = 2520
$ws.Range("K99").Value = 2520
$ws.Range("M99").Value = -1022
$ws.Range("H105").Value = 4590.6523
$ws.Range("I105").Value = 4393.0557
$ws.Range("J105").Value = 5302
$ws.Range("K105").Value = 4393.0557
$ws.Range("L105").Value = 5302
$ws.Range("M105").Value = -2646.0557
$ws.Range("N105").Value = -8796
$ws.Range("H134").Value = 944.8378
$ws.Range("I134").Value = 870.25714
$ws.Range("J134").Value = 2250
$ws.Range("K134").Value = 2610.77142
$ws.Range("L134").Value = 6750
$ws.Range("M134").Value = -75.77142000000003
$ws.Range("N134").Value = -11820
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 233.33333
$ws.Range("J22").Value = 400
$ws.Range("L22").Value = 400
$ws.Range("N22").Value = -1100
$ws.Range("H31").Value = 1542.96
$ws.Range("J31").Value = 3788.75
$ws.Range("L31").Value = 3788.75
$ws.Range("N31").Value = -4378.75
$ws.Range("H34").Value = 1542.96
$ws.Range("J34").Value = 3788.75
$ws.Range("L34").Value = 3788.75
$ws.Range("N34").Value = -4192.75
$ws.Range("H55").Value = 20000
$ws.Range("I55").Value = 20000
$ws.Range("K55").Value = 20000
$ws.Range("M55").Value = -19685
$ws.Range("H99").Value = 18050
$ws.Range("I99").Value = 20313.875
$ws.Range("K99").Value = 20313.875
$ws.Range("M99").Value = -18815.875
$ws.Range("H105").Value = 1452.8636
$ws.Range("I105").Value = 1431.4286
$ws.Range("J105").Value = 1490.375
$ws.Range("K105").Value = 1431.4286
$ws.Range("L105").Value = 1490.375
$ws.Range("M105").Value = 315.5714
$ws.Range("N105").Value = -4984.375
$ws.Range("H126").Value = 18050
$ws.Range("I126").Value = 20313.875
$ws.Range("K126").Value = 60941.625
$ws.Range("M126").Value = -58471.625
$ws.Range("H132").Value = 2940.1428
$ws.Range("I132").Value = 2464.7778
$ws.Range("J132").Value = 3795.8
$ws.Range("K132").Value = 7394.3334
$ws.Range("L132").Value = 11387.4
$ws.Range("M132").Value = -4864.3334
$ws.Range("N132").Value = -16447.4
$ws.Range("H134").Value = 2280.36
$ws.Range("I134").Value = 2353.4211
$ws.Range("K134").Value = 7060.263300000001
$ws.Range("M134").Value = -4525.263300000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6214.9473
$ws.Range("I134").Value = 2161.077
$ws.Range("J134").Value = 14998.333
$ws.Range("K134").Value = 6483.231000000001
$ws.Range("L134").Value = 44994.999
$ws.Range("M134").Value = -1413.231000000001
$ws.Range("N134").Value = -55134.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3229.1177
$ws.Range("I102").Value = 2793
$ws.Range("K102").Value = 2793
$ws.Range("M102").Value = -1171
$ws.Range("H132").Value = 2387.8965
$ws.Range("I132").Value = 2170.32
$ws.Range("J132").Value = 3747.75
$ws.Range("K132").Value = 6510.960000000001
$ws.Range("L132").Value = 11243.25
$ws.Range("M132").Value = -3980.960000000001
$ws.Range("N132").Value = -16303.25
$ws.Range("H135").Value = 172500
$ws.Range("J135").Value = 172500
$ws.Range("L135").Value = 172500
$ws.Range("N135").Value = -182640
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5158.5713
$ws.Range("I40").Value = 3542.75
$ws.Range("J40").Value = 7313
$ws.Range("K40").Value = 3542.75
$ws.Range("L40").Value = 7313
$ws.Range("M40").Value = -3406.75
$ws.Range("N40").Value = -7585
$ws.Range("H122").Value = 9330.588
$ws.Range("I122").Value = 11402.167
$ws.Range("J122").Value = 4358.8
$ws.Range("K122").Value = 34206.501
$ws.Range("L122").Value = 13076.4
$ws.Range("M122").Value = -31756.501
$ws.Range("N122").Value = -17976.4
$ws.Range("H132").Value = 519884.53
$ws.Range("I132").Value = 653854.6
$ws.Range("K132").Value = 1961563.8
$ws.Range("M132").Value = -1959033.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 99961
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 99961
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 99961
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -103275
$ws.Range("H123").Value = 40300
$ws.Range("J123").Value = 40300
$ws.Range("L123").Value = 40300
$ws.Range("N123").Value = -50100
$ws.Range("H126").Value = 2836.5334
$ws.Range("I126").Value = 2004.4546
$ws.Range("K126").Value = 6013.3638
$ws.Range("M126").Value = -3543.3638
$ws.Range("H132").Value = 52934.465
$ws.Range("I132").Value = 50286.57
$ws.Range("K132").Value = 150859.71
$ws.Range("M132").Value = -148329.71
